# Daily attendance processing - 2025-11-02 20:42:57
# Swap the order of entries in the "Recorded By" (column G) list so that
# "dnasr281@gmail.com" appears first, followed by the other recorder name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$targetEmail = "dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $value = $cell.Text

    if ($value -ne $null -and $value -like "*, $targetEmail") {
        $parts = $value -split ", ", 2
        if ($parts.Count -eq 2 -and $parts[1] -eq $targetEmail) {
            $cell.Value2 = "$($parts[1]), $($parts[0])"
        }
    }
}
